# Add Some Missing Prompts
# - Insert a new "MaxTries.wav" row into the prompts table (after NoMatch.wav)
# - Split the old combined LangMenu prompt into two separate, shorter prompts

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# 1) Insert a blank worksheet row at row 7 (pushes PhoneNumber.wav.. down by one)
$ws.Rows.Item(7).Insert()

# 2) Grow the table so it covers the newly inserted row as well
$lo.Resize($ws.Range("A1:C22"))

# 3) Fill in the new row with the "MaxTries" prompt
$ws.Range("A7").Value = "MaxTries.wav "
$ws.Range("B7").Value = "Sorry, you have exceeded the maximum number of attempts."
$ws.Range("C7").Value = "عذراً، لقد تجاوزت الحد الأقصى لعدد المحاولات. "

# 4) Split the LangMenu prompt (row 3) into two shorter, separate prompts
$ws.Range("B3").Value = "For English press 2"
$ws.Range("C3").Value = "للغة العربية اضغط 1"

# 5) Move the selection, matching the author's final cursor position
$ws.Range("A2").Select()
